$wb = $excel.ActiveWorkbook

# Rename the two shared strings that had a stray leading space so the
# shared-string table is rebuilt: " SweepingInterval" -> "SweepingInterval"
# and " SweepingFractionAvailable" -> "SweepingFractionAvailable".
$landuses = $wb.Worksheets.Item("LANDUSES")
$landuses.Range("B1").Value = "SweepingInterval"
$landuses.Range("C1").Value = "SweepingFractionAvailable"

# Move the active tab from LOADINGS to LANDUSES, with C1 selected there.
$landuses.Activate()
$landuses.Range("C1").Select()
